$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 3: second test case "Tc_02_RegisterNewUser" ---
$ws.Range("A3").Value = "Tc_02_RegisterNewUser"
$ws.Range("B3").Value = 1

# C3/D3 reuse the same "bordered data cell" formatting already used by C2/D2
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("C3").Value = "chrome"
$ws.Range("D3").Value = "http://automationpractice.com/index.php"

# --- New column E: "EmailId" header, styled like the other header cells ---
$ws.Range("E1").Value = "EmailId"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Borders.Item(9).LineStyle = -4142
$ws.Range("E1").Borders.Item(8).LineStyle = -4142

# E3 holds the new user's e-mail address as a live mailto hyperlink
$ws.Range("E3").Value = "abc@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:abc@gmail.com") | Out-Null

# Column A needs to widen to fit the longer test-case name
$ws.Columns("A").ColumnWidth = 24.42578125

$ws.Range("E3").Select()

Write-Output "done"
